$d = $word.ActiveDocument

# The last paragraph in the document is an (empty) "Prrafodelista" bullet
# item that only holds the _GoBack bookmark. We turn it into a "Hash
# tables" Heading-3 paragraph, and push the bookmark (still attached to
# the tail of the original paragraph) down into a brand-new plain
# paragraph that carries the explanatory sentence.

$n = $d.Paragraphs.Count
$target = $d.Paragraphs.Item($n)

# Split: insert a new empty paragraph right before the bookmark
# paragraph. The bookmark paragraph (and its bookmark) ends up as the
# next paragraph, unaffected.
[void]$target.Range.InsertParagraphBefore()

$heading = $d.Paragraphs.Item($n)
$body = $d.Paragraphs.Item($n + 1)

# Turn the new empty paragraph into the "Hash tables" heading. The
# pStyle has to travel inside the inserted XML itself: InsertXML on a
# whole (collapsed, paragraph-mark-only) paragraph range rewrites the
# paragraph's pPr too, so a `.Style =` done beforehand would otherwise
# be discarded.
$headingXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Ttulo3"/></w:pPr><w:r><w:t xml:space="preserve">Hash </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>tables</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$heading.Range.InsertXML($headingXml)

# Turn the (still bookmarked) paragraph into a regular paragraph
# *before* inserting the sentence: inserting into a collapsed range in
# the middle of a paragraph's content leaves the existing pPr alone, so
# the style has to be applied directly to the paragraph itself.
$body.Style = "Normal"
$bodyStart = $d.Range($body.Range.Start, $body.Range.Start)
$bodyXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Tambi&#233;n llamados hash </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>maps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>maps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>unordered</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>maps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dictionaries</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>objects</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, etc.</w:t></w:r><w:r><w:t>, guardan en memoria un valor en base a una llave dada</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$bodyStart.InsertXML($bodyXml)
